$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.153.60'
$ws.Range("E2").Value = '  -3.85%  '

$ws.Range("D3").Value = '2.882.25'
$ws.Range("E3").Value = '  -4.29%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = "'523.74"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'140.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.89%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("E8").Value = '  -3.71%  '

$ws.Range("D9").Value = '2.880.81'
$ws.Range("E9").Value = '  -4.45%  '

$ws.Range("E10").Value = '  -5.99%  '

$ws.Range("D11").Value = "'5.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.36%  '

$ws.Range("E12").Value = '  -3.64%  '

$ws.Range("D13").Value = '3.394.52'
$ws.Range("E13").Value = '  -3.99%  '

$ws.Range("E14").Value = '  +0.96%  '

$ws.Range("D15").Value = '60.328.98'
$ws.Range("E15").Value = '  -3.72%  '

$ws.Range("D16").Value = "'22.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.21%  '

$ws.Range("D17").Value = '2.891.63'
$ws.Range("E17").Value = '  -3.94%  '

$ws.Range("E18").Value = '  -6.16%  '

$ws.Range("E19").Value = '  -4.15%  '

$ws.Range("D20").Value = "'11.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.36%  '

$ws.Range("D21").Value = "'358.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -9.06%  '

$ws.Range("D22").Value = "'6.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.42%  '

$ws.Range("E23").Value = '  -0.15%  '

$ws.Range("D24").Value = "'62.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.06%  '

$ws.Range("D25").Value = '3.020.96'
$ws.Range("E25").Value = '  -3.72%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = "'0.181"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.81%  '

$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").Value = "'0.444"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.31%  '

$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("D29").Value = "'7.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.41%  '

$ws.Range("D30").Value = '0.0₃0843'
$ws.Range("E30").Value = '  -12.82%  '

$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("E32").Value = '  -5.04%  '

$ws.Range("E33").Value = '  -5.80%  '

$ws.Range("D34").Value = "'148.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.25%  '

$ws.Range("E35").Value = '  -8.79%  '

$ws.Range("D36").Value = "'5.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.53%  '

$ws.Range("D37").Value = "'0.984"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -9.01%  '

$ws.Range("E38").Value = '  -8.48%  '

$ws.Range("D39").Value = "'37.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.08%  '

$ws.Range("E40").Value = '  -6.38%  '

$ws.Range("D41").Value = '2.313.14'
$ws.Range("E41").Value = '  -5.76%  '

$ws.Range("D42").Value = "'0.639"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.63%  '

$ws.Range("E43").Value = '  -8.06%  '

$ws.Range("D44").Value = "'20.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -9.11%  '

$ws.Range("D45").Value = "'0.0565"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.12%  '

$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.14%  '

$ws.Range("E47").Value = '  +0.76%  '

$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").Value = "'10.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.52%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = "'0.0232"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.22%  '

$ws.Range("E50").Value = '  -3.42%  '

$ws.Range("D51").Value = "'247.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.42%  '
